# Populate the "Login" sheet with the expanded dashboard/menu table and
# mark it as the active/selected tab (the commit's "Additional scenario
# Dashboard executed" data dump).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# Touch the cells that stay blank in the new A1:H6 block so they are
# materialised (with the sheet's normal style) instead of left absent,
# matching the full 6x8 grid written by the dashboard export.
$blankCells = @("C2","A3","B3","C3","A4","B4","C4","D4","A5","B5","C5","D5","H5","A6","B6","C6","D6","F6","H6")
foreach ($addr in $blankCells) {
  $ws.Range($addr).Style = "Normal"
}

# --- Row 1: headers -------------------------------------------------
$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "password"
$ws.Range("C1").Value = "message"
$ws.Range("D1").Value = "title"
$ws.Range("E1").Value = "menu"
$ws.Range("F1").Value = "icon"
$ws.Range("G1").Value = "table"
$ws.Range("H1").Value = "page"

# --- Row 2 ------------------------------------------------------------
$ws.Range("A2").Value = "sdetnumpyninja@gmail.com"
$ws.Range("B2").Value = "Feb@2025"
$ws.Range("D2").Value = "LMS - Learning Management System"
$ws.Range("E2").Value = "Home"
$ws.Range("F2").Value = "User"
$ws.Range("G2").Value = "Staff Data"
$ws.Range("H2").Value = "Manage User"

# --- Row 3 --------------------------------------------------------------
$ws.Range("D3").Value = "LMS"
$ws.Range("E3").Value = "Program"
$ws.Range("F3").Value = "Staff"
$ws.Range("G3").Value = "#"
$ws.Range("H3").Value = "Manage Batch"

# --- Row 4 --------------------------------------------------------------
$ws.Range("E4").Value = "Batch"
$ws.Range("F4").Value = "Batches"
$ws.Range("G4").Value = "First Name"
$ws.Range("H4").Value = "Manage Program"

# --- Row 5 --------------------------------------------------------------
$ws.Range("E5").Value = "Class"
$ws.Range("F5").Value = "Programs"
$ws.Range("G5").Value = "Last Name"

# --- Row 6 --------------------------------------------------------------
$ws.Range("E6").Value = "Logout"
$ws.Range("G6").Value = "Phone"

# Select the Login sheet as the active tab.
$ws.Activate()
$ws.Select()
